# Update number of cities (population figures), re-sort the table by
# Population descending, then turn on AutoFilter for the data range
# (which in turn defines the hidden _FilterDatabase name, as Excel does).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update populations that changed -------------------------------------
# Buenos Aires (currently row 4, before the re-sort below)
$ws.Range("C4").Value = 2891000
# Toronto (currently row 12, before the re-sort below)
$ws.Range("C12").Value = 2800000

# --- Sort the data range A1:D13 by Population (column C) descending ------
$dataRange = $ws.Range("A1:D13")
$sortKey = $ws.Range("C1")
$dataRange.Sort($sortKey, 2)

# --- Turn on AutoFilter over the sorted table -----------------------------
[void]$ws.Range("A1:D13").AutoFilter()

# --- Mirror Excel's hidden _FilterDatabase name for the AutoFilter range --
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$D`$13")
$filterName.Visible = $false
